$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Lesson 21 "Exceptions" video uploaded to YouTube; also fill in the two
# lessons that follow it (File I/O & Date-Time API, Practice - recap) plus
# the already-known "OOP #5" / "OOP #7" hyperlinks that were still pending.
# ---------------------------------------------------------------------------

# --- Row 24: lesson 21 "Exceptions" ---
$ws.Range("C24").Value = "Exceptions"
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 44173
$ws.Range("F24").Value = "https://youtu.be/fVmROnih-Io"

# --- Row 25: lesson 22 "File I/O & Date-Time API" ---
$ws.Range("C25").Value = "File I/O & Date-Time API"
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 44177
$ws.Range("E25").NumberFormat = "d-mmm-yy"
$ws.Range("E25").Font.Size = 10

# --- Row 26: lesson 23 "Practice - recap" (reuses existing lesson text) ---
$ws.Range("C26").Value = "Practice – recap "
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 44179
$ws.Range("E26").NumberFormat = "d-mmm-yy"
$ws.Range("E26").Font.Size = 10

# Row heights: 21 & 23 shrink slightly, 24-26 grow now that they hold data
$ws.Rows.Item(21).RowHeight = 13.85
$ws.Rows.Item(23).RowHeight = 13.85
$ws.Rows.Item(24).RowHeight = 14.9
$ws.Rows.Item(25).RowHeight = 14.9
$ws.Rows.Item(26).RowHeight = 14.9

# Column F narrows now that the long mega.nz links aren't the widest entries
$ws.Columns.Item(6).ColumnWidth = 51.6

# ---------------------------------------------------------------------------
# Hyperlinks: rebuild the whole collection in row order so relationship ids
# renumber the same way Excel would after adding the new links for F21,
# F23 and F24 (F22's pre-existing link simply shifts down to rId11).
# ---------------------------------------------------------------------------
$ws.Range("F2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://youtu.be/3Q7s1cpByuk", "", "", "https://youtu.be/3Q7s1cpByuk")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://youtu.be/9fIwJtlF_Dg", "", "", "https://youtu.be/9fIwJtlF_Dg")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://youtu.be/tqJZedXiqeI", "", "", "https://youtu.be/tqJZedXiqeI")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://youtu.be/yNGEqPdB944", "", "", "https://youtu.be/yNGEqPdB944")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://youtu.be/OgDodkukz9U", "", "", "https://youtu.be/OgDodkukz9U")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://youtu.be/G0B2xyAF3RY", "", "", "https://youtu.be/G0B2xyAF3RY")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://youtu.be/GAIGLm6nZVI", "", "", "https://youtu.be/GAIGLm6nZVI")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://mega.nz/folder/99JXxCJB", "5W1WODw_dXbtaLwv7hdQSQ", "", "https://mega.nz/folder/99JXxCJB#5W1WODw_dXbtaLwv7hdQSQ")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://mega.nz/folder/99JXxCJB", "5W1WODw_dXbtaLwv7hdQSQ", "", "https://mega.nz/folder/99JXxCJB#5W1WODw_dXbtaLwv7hdQSQ")
$ws.Hyperlinks.Add($ws.Range("F21"), "https://youtu.be/1XCeWEAcA4I", "", "", "https://youtu.be/1XCeWEAcA4I")
$ws.Hyperlinks.Add($ws.Range("F22"), "https://youtu.be/WZK0AT6SJfk", "", "", "https://youtu.be/WZK0AT6SJfk")
$ws.Hyperlinks.Add($ws.Range("F23"), "https://youtu.be/XBw5dgHg3dE", "", "", "https://youtu.be/XBw5dgHg3dE")
$ws.Hyperlinks.Add($ws.Range("F24"), "https://youtu.be/fVmROnih-Io", "", "", "https://youtu.be/fVmROnih-Io")

# Scroll / selection, matching where the editor ended up after the edit
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F35").Select()

Write-Host "lesson 21 video record uploaded into youtube"
